# Apply polling place list updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (#21 投票区) - name change only
$ws.Range("E22").Value = "旧春日幼稚園"
$ws.Range("F22").Value = "キュウカスガヨウチエン"

# Row 30 (#29 投票区) - coordinates, name, kana, address change
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "34.29253484447429"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "134.05250434211888"
$ws.Range("E30").Value = "多肥小学校体育館"
$ws.Range("F30").Value = "タヒショウガッコウタイイクカン"
$ws.Range("G30").Value = "香川県高松市多肥上町902-2"

# Row 33 (#32 投票区) - coordinates, name, kana, address change
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "34.28638923532128"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "134.02356112267915"
$ws.Range("E33").Value = "高松南高校校舎棟西館"
$ws.Range("F33").Value = "タカマツミナミコウコウコウシャトウニシカン"
$ws.Range("G33").Value = "香川県高松市一宮町531"

# Row 43 (#42 投票区) - coordinates, name, kana, address change
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "34.27696"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "134.08578"
$ws.Range("E43").Value = "山田総合センター内大ホール"
$ws.Range("F43").Value = "ヤマダソウゴウセンターナイダイホール"
$ws.Range("G43").Value = "香川県高松市川島本町191-10"

# Row 65 (#65 投票区) - coordinates, name, kana, address change
$ws.Range("B65").NumberFormat = "@"
$ws.Range("B65").Value = "34.3840175"
$ws.Range("C65").NumberFormat = "@"
$ws.Range("C65").Value = "134.12656556"
$ws.Range("E65").Value = "庵治支所"
$ws.Range("F65").Value = "アジシショ"
$ws.Range("G65").Value = "香川県高松市庵治町6393-5"
